$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $savedStyle = $rng.Style
    $rng.Value = "'" + $value
    $rng.Style = $savedStyle
}

Set-TextValue "D2" "69.175.19"
Set-TextValue "E2" "  +0.10%  "
Set-TextValue "D3" "3.746.89"
Set-TextValue "E3" "  +0.41%  "
Set-TextValue "E4" "  +0.08%  "
Set-TextValue "D5" "601.72"
Set-TextValue "E5" "  +0.00%  "
Set-TextValue "D6" "166.70"
Set-TextValue "E6" "  -1.52%  "
Set-TextValue "D7" "3.748.45"
Set-TextValue "E7" "  +0.35%  "
Set-TextValue "E8" "  -0.02%  "
Set-TextValue "E9" "  +1.13%  "
Set-TextValue "D10" "0.170"
Set-TextValue "E10" "  +2.17%  "
Set-TextValue "D11" "6.41"
Set-TextValue "E11" "  +1.57%  "
Set-TextValue "E12" "  -0.52%  "
Set-TextValue "D13" "37.96"
Set-TextValue "E13" "  -1.34%  "
Set-TextValue "D14" "0.0000249"
Set-TextValue "E14" "  +1.30%  "
Set-TextValue "D15" "4.373.98"
Set-TextValue "E15" "  +0.50%  "
Set-TextValue "D16" "3.725.78"
Set-TextValue "E16" "  +0.18%  "
Set-TextValue "D17" "69.131.74"
Set-TextValue "E17" "  +0.20%  "
Set-TextValue "D18" "7.36"
Set-TextValue "E18" "  +1.02%  "
Set-TextValue "D19" "17.39"
Set-TextValue "E19" "  -0.57%  "
Set-TextValue "E20" "  -1.52%  "
Set-TextValue "D21" "11.06"
Set-TextValue "E21" "  +17.91%  "
Set-TextValue "D22" "494.54"
Set-TextValue "E22" "  -0.76%  "
Set-TextValue "E23" "  -0.04%  "
Set-TextValue "D24" "0.0000153"
Set-TextValue "E24" "  +7.72%  "
Set-TextValue "D25" "84.86"
Set-TextValue "E25" "  -0.09%  "
Set-TextValue "E26" "  -1.10%  "
Set-TextValue "D27" "12.30"
Set-TextValue "E27" "  -0.18%  "
Set-TextValue "D28" "10.10"
Set-TextValue "E28" "  -0.27%  "
Set-TextValue "E29" "  -0.18%  "
Set-TextValue "E30" "  +1.68%  "
Set-TextValue "E31" "  +2.30%  "
Set-TextValue "D32" "8.05"
Set-TextValue "E32" "  +0.74%  "
Set-TextValue "D33" "31.57"
Set-TextValue "E33" "  -0.41%  "
Set-TextValue "D34" "3.890.65"
Set-TextValue "E34" "  +0.64%  "
Set-TextValue "B35" "RenzoRestakedETH"
Set-TextValue "C35" "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
Set-TextValue "D35" "3.684.72"
Set-TextValue "E35" "  +0.41%  "
Set-TextValue "B36" "Hedera"
Set-TextValue "C36" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D36" "0.108"
Set-TextValue "E36" "  -0.58%  "
Set-TextValue "E37" "  +0.09%  "
Set-TextValue "E38" "  +0.62%  "
Set-TextValue "D39" "5.95"
Set-TextValue "E39" "  +2.21%  "
Set-TextValue "E40" "  +0.68%  "
Set-TextValue "E41" "  -0.40%  "
Set-TextValue "E42" "  +4.08%  "
Set-TextValue "D43" "430.84"
Set-TextValue "E43" "  -1.18%  "
Set-TextValue "D44" "48.78"
Set-TextValue "E44" "  -1.10%  "
Set-TextValue "E45" "  -0.22%  "
Set-TextValue "D46" "8.48"
Set-TextValue "E46" "  +0.67%  "
Set-TextValue "E47" "  +0.02%  "
Set-TextValue "D48" "40.19"
Set-TextValue "E48" "  -1.44%  "
Set-TextValue "D49" "141.02"
Set-TextValue "E49" "  -1.35%  "
Set-TextValue "D50" "2.790.37"
Set-TextValue "E50" "  +1.49%  "
Set-TextValue "E51" "  +0.01%  "
